$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume change (E) columns with freshly
# scraped values. D/E cells store plain text (not numbers/percentages),
# matching the source feed formatting exactly.

$ws.Range("D2").Value = '66.297.59'
$ws.Range("E2").Value = '  -0.83%  '

$ws.Range("D3").Value = '3.323.14'
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'588.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.57%  '

$ws.Range("D6").Value = "'183.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("D7").Value = "'0.647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.02%  '

$ws.Range("E9").Value = '  -2.67%  '

$ws.Range("D10").Value = "'6.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.19%  '

$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("D12").Value = '3.901.76'
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("E13").Value = '  -3.46%  '

$ws.Range("D14").Value = '66.312.75'
$ws.Range("E14").Value = '  -0.86%  '

$ws.Range("E15").Value = '  -3.45%  '

$ws.Range("D16").Value = "'0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.65%  '

$ws.Range("D17").Value = '3.270.81'
$ws.Range("E17").Value = '  -1.96%  '

$ws.Range("D18").Value = "'427.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.09%  '

$ws.Range("E19").Value = '  -2.71%  '

$ws.Range("D20").Value = "'13.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.11%  '

$ws.Range("E21").Value = '  -2.84%  '

$ws.Range("D22").Value = "'71.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.68%  '

$ws.Range("E23").Value = '  +0.17%  '

$ws.Range("E24").Value = '  +1.07%  '

$ws.Range("D25").Value = '3.462.54'
$ws.Range("E25").Value = '  -0.88%  '

$ws.Range("D26").Value = "'0.516"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.65%  '

$ws.Range("D27").Value = "'0.205"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.58%  '

$ws.Range("E28").Value = '  -3.61%  '

$ws.Range("D29").Value = "'8.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.89%  '

$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("E31").Value = '  -0.65%  '

$ws.Range("D32").Value = "'22.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.16%  '

$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").Value = "'5.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.36%  '

$ws.Range("E35").Value = '  -3.07%  '

$ws.Range("E36").Value = '  -3.76%  '

$ws.Range("D37").Value = "'159.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.96%  '

$ws.Range("D38").Value = "'1.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.87%  '

$ws.Range("D39").Value = '2.886.06'
$ws.Range("E39").Value = '  +1.61%  '

$ws.Range("E40").Value = '  -2.39%  '

$ws.Range("D41").Value = "'26.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.97%  '

$ws.Range("E42").Value = '  -3.31%  '

$ws.Range("D43").Value = "'4.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.70%  '

$ws.Range("D44").Value = "'40.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("D45").Value = "'0.0666"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.01%  '

$ws.Range("D46").Value = "'5.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.19%  '

$ws.Range("E47").Value = '  -2.62%  '

$ws.Range("D48").Value = "'23.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.54%  '

$ws.Range("D49").Value = "'314.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.00%  '

$ws.Range("E50").Value = '  -0.37%  '

$ws.Range("E51").Value = '  +5.10%  '
